$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repulled data / mean calculation fix
$ws.Range("F5").Value = 5
$ws.Range("F8").Value = 6
$ws.Range("F13").Value = -2
$ws.Range("F19").Value = -2
$ws.Range("F34").Value = -8
$ws.Range("F36").Value = -4
$ws.Range("F37").Value = 2
$ws.Range("F44").Value = -6
$ws.Range("F51").Value = 2
$ws.Range("F52").Value = 1
$ws.Range("F60").Value = -6
